$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.079.43"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "3.428.81"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'579.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'145.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("D8").Value = "'0.475"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'7.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "'0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "4.013.96"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'28.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "3.422.41"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "62.107.23"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'14.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").Value = "'9.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").Value = "'392.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").Value = "'74.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'0.0000116"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "3.566.25"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").Value = "'0.186"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "'7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").Value = "'2.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'23.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("E35").Value = "  +5.93%  "
$ws.Range("D36").Value = "'7.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "'167.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").Value = "'1.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("D39").Value = "3.461.05"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "'28.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.76%  "
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("D46").Value = "2.513.54"
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("D47").Value = "'23.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "'6.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0265"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.85%  "
